$p = $ppt.ActivePresentation
$d = $p.Designs.Item(1)
Write-Host "before:" $d.Name
try {
    $d.Name = "Office Theme"
    Write-Host "set OK, now:" $d.Name
} catch {
    Write-Host "set failed:" $_
}
